# "Todas las bases homologadas para unirse"
#
# Homologate this sheet to the common municipio-catalog key (CVE_MUN) so it
# can be merged/joined with the other standardized bases:
#   1. Column A header "Municipio" -> "CVE_MUN"
#   2. Column A values: municipio name -> its INEGI municipality code (CVE_MUN)
#   3. Drop the stray "Publico a/" catalog row (B=0, C=0.09) that isn't a
#      real municipio and has no CVE_MUN code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header -------------------------------------------------------------
$ws.Range("A1").Value = "CVE_MUN"

# --- 2) Municipio name -> CVE_MUN code, by original row number -------------
$cveMap = @{
    2  = "13003"
    3  = "13007"
    4  = "13008"
    5  = "13010"
    6  = "13011"
    7  = "13013"
    8  = "13012"
    9  = "13014"
    10 = "13015"
    11 = "13017"
    12 = "13018"
    13 = "13019"
    14 = "13016"
    15 = "13009"
    16 = "13021"
    17 = "13022"
    18 = "13024"
    19 = "13025"
    20 = "13027"
    21 = "13029"
    22 = "13030"
    23 = "13031"
    24 = "13034"
    25 = "13051"
    26 = "13038"
    27 = "13042"
    28 = "13044"
    30 = "13048"
    31 = "13049"
    32 = "13052"
    33 = "13053"
    34 = "13046"
    35 = "13054"
    36 = "13056"
    37 = "13058"
    38 = "13059"
    39 = "13060"
    40 = "13061"
    41 = "13062"
    42 = "13063"
    43 = "13069"
    44 = "13073"
    45 = "13076"
    46 = "13077"
    47 = "13066"
    48 = "13078"
    49 = "13080"
    50 = "13082"
    51 = "13084"
}

foreach ($row in $cveMap.Keys) {
    $ws.Cells.Item($row, 1).Value = $cveMap[$row]
}

# --- 3) Remove the "Publico a/" row (old row 29) ----------------------------
$ws.Rows(29).Delete()
